$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47
$prev = $row - 1

# Copy formatting from the row above so new row matches existing style pattern
$ws.Range("A$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E$prev").Copy()
$ws.Range("E$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 46
$ws.Cells.Item($row, 2).Value = "portugal"
$ws.Cells.Item($row, 3).Value = "liga-portugal"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45191.88541666666
$ws.Cells.Item($row, 6).Value = "Famalicao"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Arouca"
$ws.Cells.Item($row, 9).Value = 0

$ws.Cells.Item($row, 10).Value = 2.02
$ws.Cells.Item($row, 11).Value = "17/09/2023 15:42"
$ws.Cells.Item($row, 12).Value = 2.14
$ws.Cells.Item($row, 13).Value = "22/09/2023 21:00"

$ws.Cells.Item($row, 14).Value = 3.56
$ws.Cells.Item($row, 15).Value = "17/09/2023 15:42"
$ws.Cells.Item($row, 16).Value = 3.4
$ws.Cells.Item($row, 17).Value = "22/09/2023 20:58"

$ws.Cells.Item($row, 18).Value = 3.89
$ws.Cells.Item($row, 19).Value = "17/09/2023 15:42"
$ws.Cells.Item($row, 20).Value = 3.74
$ws.Cells.Item($row, 21).Value = "22/09/2023 21:00"

$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/famalicao-arouca/rsDhXiL9/"
